# Update the "Mobile networks" reference list: remove "O2" as a mobile
# network option. This is done by deleting the entire row that holds it
# (row 3 on the hidden "Mobile networks" sheet) rather than just clearing
# the cell, so that the table/list below it shifts up, the named table
# shrinks, and the now-unused "O2" shared string drops out of the workbook
# on save - exactly what Excel does for a real "delete row" edit.

$wb = $excel.ActiveWorkbook

$networksSheet = $wb.Worksheets.Item("Mobile networks")
$mainSheet = $wb.Worksheets.Item("Extra mobile data requests")

# Row 3 is "O2" (row 1 is the header "Mobile networks", row 2 is "EE").
$networksSheet.Rows.Item(3).Delete() | Out-Null

# Mirror the cursor position Excel leaves behind after deleting a row:
# the row that slid up into the deleted row's place ends up selected.
$networksSheet.Rows.Item(3).Select() | Out-Null

# Restore the main (visible) sheet as the active tab/selection, matching
# where the author's cursor was left in the saved workbook.
$mainSheet.Activate() | Out-Null
$mainSheet.Range("C2").Select() | Out-Null
